# "changes to template and ppt"
#
# 1) Every "datetimeFigureOut" date field (slide master + all 11 slide
#    layouts) moves from 1/7/2014 to 1/8/2014.
# 2) Slide 2 ("Once upon a time") content placeholder text is expanded:
#      - paragraph 1 gets a new trailing sentence
#      - paragraph 2 gets an extra word ("theoretical") inserted
# 3) Slide 3's (empty) title becomes "Callbacks".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder fields: slide master + all custom (slide) layouts
# ---------------------------------------------------------------------
$sm = $p.SlideMaster
$sm.Shapes.Item(3).TextFrame.TextRange.Text = "1/8/2014"

$layouts = $sm.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "1/8/2014"
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2 body copy
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# --- paragraph 1: "In the mythical land of JS{ellipsis}." ---------------
$para1 = $tr.Paragraphs(1, 1)
$tail = $tr.Characters($para1.Start + 26, 2)
$tail.Text = [char]0x2026 + ". there was something called a callback"

# --- paragraph 2: the quoted sentence -----------------------------------
$para2 = $tr.Paragraphs(2, 1)
$p2start = $para2.Start

$seg1 = $tr.Characters($p2start, 5)
$seg1.Text = [char]0x201C + "Cut "

$seg2 = $tr.Characters($p2start + 5, 4)
$seg2.Text = "the "

$seg3 = $tr.Characters($p2start + 9, 5)
$seg3.Text = "theoretical crap "

# ---------------------------------------------------------------------
# 3) Slide 3 title
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Callbacks"
